$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers must be forced to
# text (matching the original inlineStr cell type) to avoid Excel
# auto-converting them to floating point numbers (precision loss).
$numericLookingCells = @(
    ,@(5, 4)
    ,@(6, 4)
    ,@(8, 4)
    ,@(10, 4)
    ,@(11, 4)
    ,@(12, 4)
    ,@(16, 4)
    ,@(20, 4)
    ,@(21, 4)
    ,@(22, 4)
    ,@(23, 4)
    ,@(24, 4)
    ,@(25, 4)
    ,@(29, 4)
    ,@(31, 4)
    ,@(32, 4)
    ,@(33, 4)
    ,@(34, 4)
    ,@(35, 4)
    ,@(36, 4)
    ,@(37, 4)
    ,@(38, 4)
    ,@(40, 4)
    ,@(41, 4)
    ,@(42, 4)
    ,@(43, 4)
    ,@(45, 4)
    ,@(46, 4)
    ,@(47, 4)
    ,@(49, 4)
    ,@(50, 4)
)
foreach ($coord in $numericLookingCells) {
    $ws.Cells.Item($coord[0], $coord[1]).NumberFormat = "@"
}

# Apply all cell value updates
$ws.Cells.Item(2, 4).Value = '57.319.11'
$ws.Cells.Item(2, 5).Value = '  +4.25%  '
$ws.Cells.Item(3, 4).Value = '2.518.98'
$ws.Cells.Item(3, 5).Value = '  +3.20%  '
$ws.Cells.Item(4, 5).Value = '  -0.11%  '
$ws.Cells.Item(5, 4).Value = '496.78'
$ws.Cells.Item(5, 5).Value = '  +3.38%  '
$ws.Cells.Item(6, 4).Value = '154.05'
$ws.Cells.Item(6, 5).Value = '  +10.71%  '
$ws.Cells.Item(7, 5).Value = '  -0.08%  '
$ws.Cells.Item(8, 4).Value = '0.518'
$ws.Cells.Item(8, 5).Value = '  +3.20%  '
$ws.Cells.Item(9, 4).Value = '2.537.04'
$ws.Cells.Item(9, 5).Value = '  +2.79%  '
$ws.Cells.Item(10, 4).Value = '0.102'
$ws.Cells.Item(10, 5).Value = '  +5.21%  '
$ws.Cells.Item(11, 4).Value = '5.77'
$ws.Cells.Item(11, 5).Value = '  +5.62%  '
$ws.Cells.Item(12, 4).Value = '0.340'
$ws.Cells.Item(12, 5).Value = '  +4.82%  '
$ws.Cells.Item(13, 5).Value = '  +1.67%  '
$ws.Cells.Item(14, 4).Value = '2.952.86'
$ws.Cells.Item(14, 5).Value = '  +2.49%  '
$ws.Cells.Item(15, 4).Value = '57.388.49'
$ws.Cells.Item(15, 5).Value = '  +4.11%  '
$ws.Cells.Item(16, 4).Value = '21.48'
$ws.Cells.Item(16, 5).Value = '  +4.77%  '
$ws.Cells.Item(17, 5).Value = '  +3.51%  '
$ws.Cells.Item(18, 4).Value = '2.531.56'
$ws.Cells.Item(18, 5).Value = '  +2.95%  '
$ws.Cells.Item(19, 5).Value = '  +6.14%  '
$ws.Cells.Item(20, 4).Value = '10.38'
$ws.Cells.Item(20, 5).Value = '  +5.52%  '
$ws.Cells.Item(21, 4).Value = '325.12'
$ws.Cells.Item(21, 5).Value = '  +3.56%  '
$ws.Cells.Item(22, 2).Value = 'Uniswap'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(22, 4).Value = '5.98'
$ws.Cells.Item(22, 5).Value = '  +6.19%  '
$ws.Cells.Item(23, 2).Value = 'Dai'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(23, 4).Value = '0.998'
$ws.Cells.Item(23, 5).Value = '  +0.09%  '
$ws.Cells.Item(24, 4).Value = '58.67'
$ws.Cells.Item(24, 5).Value = '  +2.47%  '
$ws.Cells.Item(25, 4).Value = '0.412'
$ws.Cells.Item(25, 5).Value = '  +2.16%  '
$ws.Cells.Item(26, 5).Value = '  +2.07%  '
$ws.Cells.Item(27, 5).Value = '  -0.54%  '
$ws.Cells.Item(28, 4).Value = '2.619.73'
$ws.Cells.Item(28, 5).Value = '  +2.66%  '
$ws.Cells.Item(29, 4).Value = '7.66'
$ws.Cells.Item(29, 5).Value = '  +4.32%  '
$ws.Cells.Item(30, 4).Value = '0.0₃0834'
$ws.Cells.Item(30, 5).Value = '  +7.82%  '
$ws.Cells.Item(31, 4).Value = '0.998'
$ws.Cells.Item(31, 5).Value = '  -0.10%  '
$ws.Cells.Item(32, 4).Value = '151.93'
$ws.Cells.Item(32, 5).Value = '  +2.11%  '
$ws.Cells.Item(33, 2).Value = 'PancakeSwap'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(33, 4).Value = '1.54'
$ws.Cells.Item(33, 5).Value = '  +4.84%  '
$ws.Cells.Item(34, 2).Value = 'EthereumClassic'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(34, 4).Value = '18.46'
$ws.Cells.Item(34, 5).Value = '  +2.74%  '
$ws.Cells.Item(35, 4).Value = '5.30'
$ws.Cells.Item(35, 5).Value = '  +2.77%  '
$ws.Cells.Item(36, 2).Value = 'Fetch.AI'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(36, 4).Value = '0.916'
$ws.Cells.Item(36, 5).Value = '  +7.39%  '
$ws.Cells.Item(37, 2).Value = 'NEARProtocol'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(37, 4).Value = '3.86'
$ws.Cells.Item(37, 5).Value = '  +7.16%  '
$ws.Cells.Item(38, 4).Value = '1.17'
$ws.Cells.Item(38, 5).Value = '  +4.84%  '
$ws.Cells.Item(39, 5).Value = '  +10.38%  '
$ws.Cells.Item(40, 4).Value = '34.44'
$ws.Cells.Item(40, 5).Value = '  +3.32%  '
$ws.Cells.Item(41, 4).Value = '3.56'
$ws.Cells.Item(41, 5).Value = '  +4.53%  '
$ws.Cells.Item(42, 4).Value = '0.621'
$ws.Cells.Item(42, 5).Value = '  +2.94%  '
$ws.Cells.Item(43, 4).Value = '0.0565'
$ws.Cells.Item(43, 5).Value = '  +4.12%  '
$ws.Cells.Item(44, 5).Value = '  +0.21%  '
$ws.Cells.Item(45, 4).Value = '4.95'
$ws.Cells.Item(45, 5).Value = '  +6.29%  '
$ws.Cells.Item(46, 4).Value = '269.94'
$ws.Cells.Item(46, 5).Value = '  +5.12%  '
$ws.Cells.Item(47, 4).Value = '0.0953'
$ws.Cells.Item(47, 5).Value = '  +6.92%  '
$ws.Cells.Item(48, 5).Value = '  +4.22%  '
$ws.Cells.Item(49, 4).Value = '10.22'
$ws.Cells.Item(49, 5).Value = '  +0.89%  '
$ws.Cells.Item(50, 4).Value = '18.13'
$ws.Cells.Item(50, 5).Value = '  +6.12%  '
$ws.Cells.Item(51, 4).Value = '1.907.09'
$ws.Cells.Item(51, 5).Value = '  -1.20%  '

# Restore default (unstyled) cell style for the forced-text cells so
# the saved workbook does not carry a stray explicit style index.
foreach ($coord in $numericLookingCells) {
    $ws.Cells.Item($coord[0], $coord[1]).Style = "Normal"
}
